$wb = $excel.ActiveWorkbook

function Set-Cell {
    param(
        [string]$SheetName,
        [string]$Cell,
        [double]$Value
    )
    $ws = $wb.Worksheets($SheetName)
    $ws.Range($Cell).Value = $Value
}

# ALC sheet
Set-Cell "ALC" "H17" 1528.8096
Set-Cell "ALC" "J17" 1547.75
Set-Cell "ALC" "L17" 4643.25
Set-Cell "ALC" "N17" -4979.25

Set-Cell "ALC" "H40" 3067.25
Set-Cell "ALC" "I40" 1898
Set-Cell "ALC" "J40" 3457
Set-Cell "ALC" "K40" 1898
Set-Cell "ALC" "L40" 3457
Set-Cell "ALC" "M40" -1723
Set-Cell "ALC" "N40" -3807

Set-Cell "ALC" "H42" 156.83333
Set-Cell "ALC" "I42" 108.15385
Set-Cell "ALC" "K42" 324.46155
Set-Cell "ALC" "M42" -94.46154999999999

Set-Cell "ALC" "H51" 4766.3335
Set-Cell "ALC" "I51" 3889.3333
Set-Cell "ALC" "J51" 7397.3335
Set-Cell "ALC" "K51" 3889.3333
Set-Cell "ALC" "L51" 7397.3335
Set-Cell "ALC" "M51" -3405.3333
Set-Cell "ALC" "N51" -8365.333500000001

Set-Cell "ALC" "H55" 184.9
Set-Cell "ALC" "I55" 150.16667
Set-Cell "ALC" "J55" 237
Set-Cell "ALC" "K55" 150.16667
Set-Cell "ALC" "L55" 237
Set-Cell "ALC" "M55" 63.83332999999999
Set-Cell "ALC" "N55" -665

Set-Cell "ALC" "H69" 58859930
Set-Cell "ALC" "J69" 62507424
Set-Cell "ALC" "L69" 187522272
Set-Cell "ALC" "N69" -187524020

Set-Cell "ALC" "H70" 2720.4167
Set-Cell "ALC" "I70" 1589.6
Set-Cell "ALC" "J70" 3528.1428
Set-Cell "ALC" "K70" 4768.799999999999
Set-Cell "ALC" "L70" 10584.4284
Set-Cell "ALC" "M70" -4498.799999999999
Set-Cell "ALC" "N70" -11124.4284

Set-Cell "ALC" "H72" 58859930
Set-Cell "ALC" "J72" 62507424
Set-Cell "ALC" "L72" 562566816
Set-Cell "ALC" "N72" -562575552

Set-Cell "ALC" "H73" 2720.4167
Set-Cell "ALC" "I73" 1589.6
Set-Cell "ALC" "J73" 3528.1428
Set-Cell "ALC" "K73" 4768.799999999999
Set-Cell "ALC" "L73" 10584.4284
Set-Cell "ALC" "M73" -3832.799999999999
Set-Cell "ALC" "N73" -12456.4284

Set-Cell "ALC" "H100" 2780
Set-Cell "ALC" "I100" 2780
Set-Cell "ALC" "K100" 2780
Set-Cell "ALC" "M100" -2239

Set-Cell "ALC" "H129" 33334394
Set-Cell "ALC" "I129" 41667644
Set-Cell "ALC" "K129" 125002932
Set-Cell "ALC" "M129" -124997932

Set-Cell "ALC" "H131" 4681.25
Set-Cell "ALC" "J131" 5895
Set-Cell "ALC" "L131" 17685
Set-Cell "ALC" "N131" -27765

Set-Cell "ALC" "H138" 3156.7832
Set-Cell "ALC" "I138" 2759.7334
Set-Cell "ALC" "J138" 3244.3677
Set-Cell "ALC" "K138" 8279.200199999999
Set-Cell "ALC" "L138" 9733.1031
Set-Cell "ALC" "M138" -3139.200199999999
Set-Cell "ALC" "N138" -20013.1031

Set-Cell "ALC" "H141" 2373.5454
Set-Cell "ALC" "I141" 2069.875
Set-Cell "ALC" "K141" 6209.625
Set-Cell "ALC" "M141" -1029.625

# ARM sheet
Set-Cell "ARM" "H32" 5037.143
Set-Cell "ARM" "I32" 3118.2104
Set-Cell "ARM" "K32" 3118.2104
Set-Cell "ARM" "M32" -2831.2104

# BSM sheet
Set-Cell "BSM" "H20" 1691.9231
Set-Cell "BSM" "I20" 1486.4667
Set-Cell "BSM" "J20" 1972.091
Set-Cell "BSM" "K20" 1486.4667
Set-Cell "BSM" "L20" 1972.091
Set-Cell "BSM" "M20" -1239.4667
Set-Cell "BSM" "N20" -2466.091

# CRP sheet
Set-Cell "CRP" "H31" 22420.637
Set-Cell "CRP" "I31" 1300.4286
Set-Cell "CRP" "J31" 59381
Set-Cell "CRP" "K31" 1300.4286
Set-Cell "CRP" "L31" 59381
Set-Cell "CRP" "M31" -1005.4286
Set-Cell "CRP" "N31" -59971

Set-Cell "CRP" "H34" 22420.637
Set-Cell "CRP" "I34" 1300.4286
Set-Cell "CRP" "J34" 59381
Set-Cell "CRP" "K34" 1300.4286
Set-Cell "CRP" "L34" 59381
Set-Cell "CRP" "M34" -1098.4286
Set-Cell "CRP" "N34" -59785

Set-Cell "CRP" "H69" 39039.2
Set-Cell "CRP" "I69" 16499
Set-Cell "CRP" "K69" 16499
Set-Cell "CRP" "M69" -15750

Set-Cell "CRP" "H72" 39039.2
Set-Cell "CRP" "I72" 16499
Set-Cell "CRP" "K72" 49497
Set-Cell "CRP" "M72" -45753

Set-Cell "CRP" "H97" 28493.6
Set-Cell "CRP" "J97" 28493.6
Set-Cell "CRP" "L97" 28493.6
Set-Cell "CRP" "N97" -30475.6

Set-Cell "CRP" "H109" 34694.2
Set-Cell "CRP" "J109" 39617.75
Set-Cell "CRP" "L109" 39617.75
Set-Cell "CRP" "N109" -41697.75

Set-Cell "CRP" "H132" 128577.53
Set-Cell "CRP" "I132" 92843.63
Set-Cell "CRP" "J132" 226845.75
Set-Cell "CRP" "K132" 278530.89
Set-Cell "CRP" "L132" 680537.25
Set-Cell "CRP" "M132" -276000.89
Set-Cell "CRP" "N132" -685597.25

# CUL sheet
Set-Cell "CUL" "H37" 66129.39999999999
Set-Cell "CUL" "J37" 66129.39999999999
Set-Cell "CUL" "L37" 198388.2
Set-Cell "CUL" "N37" -198612.2

Set-Cell "CUL" "H45" 2240.25
Set-Cell "CUL" "I45" 1015
Set-Cell "CUL" "J45" 3465.5
Set-Cell "CUL" "K45" 3045
Set-Cell "CUL" "L45" 10396.5
Set-Cell "CUL" "M45" -2513
Set-Cell "CUL" "N45" -11460.5

Set-Cell "CUL" "H107" 216.11765
Set-Cell "CUL" "I107" 195
Set-Cell "CUL" "J107" 374.5
Set-Cell "CUL" "K107" 585
Set-Cell "CUL" "L107" 1123.5
Set-Cell "CUL" "M107" 1335
Set-Cell "CUL" "N107" -4963.5

Set-Cell "CUL" "H116" 2724.75
Set-Cell "CUL" "I116" 1000
Set-Cell "CUL" "J116" 4449.5
Set-Cell "CUL" "K116" 3000
Set-Cell "CUL" "L116" 13348.5
Set-Cell "CUL" "M116" 442
Set-Cell "CUL" "N116" -20232.5

# GSM sheet
Set-Cell "GSM" "H70" 6672409
Set-Cell "GSM" "I70" 11117132
Set-Cell "GSM" "J70" 5323.8335
Set-Cell "GSM" "K70" 11117132
Set-Cell "GSM" "L70" 5323.8335
Set-Cell "GSM" "M70" -11116862
Set-Cell "GSM" "N70" -5863.8335

Set-Cell "GSM" "H73" 6672409
Set-Cell "GSM" "I73" 11117132
Set-Cell "GSM" "J73" 5323.8335
Set-Cell "GSM" "K73" 11117132
Set-Cell "GSM" "L73" 5323.8335
Set-Cell "GSM" "M73" -11116196
Set-Cell "GSM" "N73" -7195.8335

Set-Cell "GSM" "H107" 1126.9
Set-Cell "GSM" "I107" 1358.2
Set-Cell "GSM" "J107" 895.6
Set-Cell "GSM" "K107" 1358.2
Set-Cell "GSM" "L107" 895.6
Set-Cell "GSM" "M107" 561.8
Set-Cell "GSM" "N107" -4735.6

Set-Cell "GSM" "H113" 18520570
Set-Cell "GSM" "I113" 55556260
Set-Cell "GSM" "J113" 2723.6667
Set-Cell "GSM" "K113" 55556260
Set-Cell "GSM" "L113" 2723.6667
Set-Cell "GSM" "M113" -55554090
Set-Cell "GSM" "N113" -7063.6667

Set-Cell "GSM" "H132" 3252.8076
Set-Cell "GSM" "I132" 2815.5833
Set-Cell "GSM" "K132" 8446.749899999999
Set-Cell "GSM" "M132" -5916.749899999999

# LTW sheet
Set-Cell "LTW" "H7" 3806.0952
Set-Cell "LTW" "I7" 2304.2856
Set-Cell "LTW" "J7" 6809.7144
Set-Cell "LTW" "K7" 2304.2856
Set-Cell "LTW" "L7" 6809.7144
Set-Cell "LTW" "M7" -2192.2856
Set-Cell "LTW" "N7" -7033.7144

Set-Cell "LTW" "H40" 8159.4165
Set-Cell "LTW" "J40" 8999.5
Set-Cell "LTW" "L40" 8999.5
Set-Cell "LTW" "N40" -9271.5

Set-Cell "LTW" "H126" 3806.0952
Set-Cell "LTW" "I126" 2304.2856
Set-Cell "LTW" "J126" 6809.7144
Set-Cell "LTW" "K126" 6912.8568
Set-Cell "LTW" "L126" 20429.1432
Set-Cell "LTW" "M126" -4442.8568
Set-Cell "LTW" "N126" -25369.1432

Set-Cell "LTW" "H132" 4119.34
Set-Cell "LTW" "I132" 3314.5715
Set-Cell "LTW" "J132" 8344.375
Set-Cell "LTW" "K132" 9943.7145
Set-Cell "LTW" "L132" 25033.125
Set-Cell "LTW" "M132" -7413.7145
Set-Cell "LTW" "N132" -30093.125

# WVR sheet
Set-Cell "WVR" "H62" 8817.643
Set-Cell "WVR" "I62" 2500
Set-Cell "WVR" "J62" 9303.615
Set-Cell "WVR" "K62" 2500
Set-Cell "WVR" "L62" 9303.615
Set-Cell "WVR" "M62" -1876
Set-Cell "WVR" "N62" -10551.615

Set-Cell "WVR" "H65" 8817.643
Set-Cell "WVR" "I65" 2500
Set-Cell "WVR" "J65" 9303.615
Set-Cell "WVR" "K65" 12500
Set-Cell "WVR" "L65" 46518.075
Set-Cell "WVR" "M65" -9380
Set-Cell "WVR" "N65" -52758.075

Set-Cell "WVR" "H81" 166666670
Set-Cell "WVR" "I81" 166666670
Set-Cell "WVR" "K81" 333333340
Set-Cell "WVR" "M81" -333332279

Set-Cell "WVR" "H84" 166666670
Set-Cell "WVR" "I84" 166666670
Set-Cell "WVR" "K84" 1666666700
Set-Cell "WVR" "M84" -1666661396
